$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from Type de composan")

# Update the "Date" metadata value (row 8, column B)
$ws1.Cells.Item(8,2).Value = "2024-07-01T07:50:29+00:00"

# Insert a new "Jurisdiction" row before "Description" (currently row 11),
# shifting Description/Purpose/Copyright/Immutable down by one row.
# Shift bottom-up so we don't clobber values we still need to read.
for ($r = 14; $r -ge 11; $r--) {
    $dst = $r + 1
    $ws1.Cells.Item($dst,1).Value = $ws1.Cells.Item($r,1).Value()
    $ws1.Cells.Item($dst,2).Value = $ws1.Cells.Item($r,2).Value()
}

# Row 15 did not exist before -- give it the same body style as the rest of the table.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 11: "Jurisdiction" with an empty value.
$ws1.Cells.Item(11,1).Value = "Jurisdiction"
$ws2.Range("A3").Copy()
$ws1.Range("B11").PasteSpecial(-4163)
$excel.CutCopyMode = 0

Write-Output "done"
